$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CountryFilter")

# The B2 cell ("BelgiumTest") should read "Belgium" instead, keeping its
# existing quote-prefixed text style. Leading with an apostrophe preserves
# the cell's "quote prefix" formatting (style index with quotePrefix="1").
$ws.Range("B2").Value = "'Belgium"

# Update the sheet's remembered selection/active cell to C4.
$ws.Range("C4").Select()
